# Update cryptocurrency price/volume data per GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.745.03"
$ws.Range("E2").Value = "  -4.99%  "
$ws.Range("D3").Value = "2.214.79"
$ws.Range("E3").Value = "  -7.05%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "296.09"
$ws.Range("E5").Value = "  -5.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "79.28"
$ws.Range("E6").Value = "  -10.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.455"
$ws.Range("E9").Value = "  -7.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0770"
$ws.Range("E10").Value = "  -6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.77"
$ws.Range("E11").Value = "  -10.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.07"
$ws.Range("E12").Value = "  -13.25%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "2.565.34"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.07"
$ws.Range("E15").Value = "  -7.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.91"
$ws.Range("E16").Value = "  -8.11%  "
$ws.Range("D17").Value = "2.238.60"
$ws.Range("E17").Value = "  -6.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.710"
$ws.Range("E18").Value = "  -6.89%  "
$ws.Range("D19").Value = "38.711.71"
$ws.Range("E19").Value = "  -4.99%  "
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -6.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.71"
$ws.Range("E21").Value = "  -7.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.64"
$ws.Range("E22").Value = "  -6.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.78"
$ws.Range("E23").Value = "  -9.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "224.49"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -10.51%  "
$ws.Range("E27").Value = "  -6.44%  "
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.03"
$ws.Range("E29").Value = "  -7.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.83"
$ws.Range("E30").Value = "  -5.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "148.21"
$ws.Range("E31").Value = "  -4.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.91"
$ws.Range("E32").Value = "  -8.78%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  -8.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.31"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0680"
$ws.Range("E36").Value = "  -7.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.63"
$ws.Range("E38").Value = "  -6.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0947"
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.58"
$ws.Range("E40").Value = "  -8.59%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.29"
$ws.Range("E41").Value = "  -11.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.57"
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("D43").Value = "1.897.00"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.01"
$ws.Range("E44").Value = "  -11.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0252"
$ws.Range("E45").Value = "  -6.93%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.93"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.99"
$ws.Range("E47").Value = "  -9.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  -10.92%  "
$ws.Range("D49").Value = "2.431.98"
$ws.Range("E49").Value = "  -7.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.08"
$ws.Range("E50").Value = "  -6.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.91"
$ws.Range("E51").Value = "  -7.27%  "
